$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; unprotect so we can edit, then re-protect after.
$ws.Unprotect()

# Update the confidential disclosure date text (A11)
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-21 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) values for rows 2-8
$ws.Range("D2").Value = 0.5002334823684389
$ws.Range("E2").Value = 0.003574559924462317

$ws.Range("D3").Value = 0.2453849852051156
$ws.Range("E3").Value = -0.004552797767660488

$ws.Range("D4").Value = 0.09518834627773895
$ws.Range("E4").Value = 0.000888099467140302

$ws.Range("D5").Value = 0.1025367150065904
$ws.Range("E5").Value = 0.003832492054589576

$ws.Range("D6").Value = 0.02984951313750463
$ws.Range("E6").Value = 0.007513727001252413

$ws.Range("D7").Value = 0.02680695800461137
$ws.Range("E7").Value = 0.001730512074254786

$ws.Range("D8").Value = 0.9999999999999998
$ws.Range("E8").Value = 0.001419105068599213

# Restore sheet protection (password unknown/not reproducible; re-protect
# with same relevant options so the sheet remains protected as before).
$ws.Protect()
